$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are forced to text via a leading apostrophe (mirrors typing
# '123 into Excel) so numeric-looking strings like "1.660" or "8.200" keep their
# exact text, matching the inline-string cells in the source file.

$ws.Range("D2").Formula = "'27.509.15"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Formula = "'1.725.44"
$ws.Range("E3").Value = "  +4.86%  "
$ws.Range("D4").Formula = "'1.004"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Formula = "'225.81"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Formula = "'0.5348"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").Formula = "'1.004"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Formula = "'0.2664"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Formula = "'0.06584"
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").Formula = "'21.61"
$ws.Range("E10").Value = "  +6.90%  "
$ws.Range("D11").Formula = "'0.07650"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Formula = "'4.593"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Formula = "'1.725.12"
$ws.Range("E13").Value = "  +4.93%  "
$ws.Range("D14").Formula = "'1.963.17"
$ws.Range("E14").Value = "  +4.86%  "
$ws.Range("D15").Formula = "'0.5809"
$ws.Range("E15").Value = "  +4.52%  "
$ws.Range("D16").Formula = "'0.0₅8290"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Formula = "'67.94"
$ws.Range("E17").Value = "  +4.65%  "
$ws.Range("D18").Formula = "'27.519.97"
$ws.Range("E18").Value = "  +5.75%  "
$ws.Range("D19").Formula = "'218.97"
$ws.Range("E19").Value = "  +13.64%  "
$ws.Range("D20").Formula = "'1.004"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Formula = "'6.055"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Formula = "'142.98"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Formula = "'1.752"
$ws.Range("E26").Value = "  +16.61%  "
$ws.Range("D27").Formula = "'0.1232"
$ws.Range("D28").Formula = "'7.347"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("E29").Value = "  +4.32%  "
$ws.Range("D30").Formula = "'0.05506"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("D31").Formula = "'1.302"
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").Formula = "'3.563"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Formula = "'3.446"
$ws.Range("E33").Value = "  +4.01%  "
$ws.Range("D34").Formula = "'1.660"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").Formula = "'2.862"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").Formula = "'2.424"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Formula = "'0.5939"
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("D39").Formula = "'0.01649"
$ws.Range("E39").Value = "  +5.10%  "
$ws.Range("D40").Formula = "'5.903"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("D41").Formula = "'1.050.37"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Formula = "'0.8486"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Formula = "'101.39"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Formula = "'1.869.49"
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("D47").Formula = "'58.83"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("D48").Formula = "'0.4474"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("D49").Formula = "'8.200"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Formula = "'1.003"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Formula = "'0.05247"
$ws.Range("E51").Value = "  +2.96%  "
